# Refresh the crypto "Price" (column D) and "Volume(1h)" (column E) snapshot
# values, as produced by the scheduled GitHub Actions scraper job.
#
# These columns hold their numbers as literal text (trailing zeros like
# "41.60" and percent signs like "-0.34%" must survive byte-for-byte), so a
# plain `.Value = "..."` assignment is not safe here: Excel's COM layer
# auto-converts numeric-/percent-looking strings into real numbers the same
# way typing them into a cell would. Prefixing with a single quote forces
# Excel's "quoted text" entry mode so the text is stored verbatim.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'330.75"
$ws.Range("E2").Value = "'-0.34%"

# Row 3
$ws.Range("D3").Value = "'41.60"
$ws.Range("E3").Value = "'0.83%"

# Row 4
$ws.Range("D4").Value = "'5.683"
$ws.Range("E4").Value = "'-0.09%"

# Row 5
$ws.Range("D5").Value = "'0.08393"
$ws.Range("E5").Value = "'4.00%"

# Row 6
$ws.Range("D6").Value = "'8.810"
$ws.Range("E6").Value = "'0.74%"

# Row 7
$ws.Range("D7").Value = "'2.008"
$ws.Range("E7").Value = "'-1.47%"

# Row 8
$ws.Range("D8").Value = "'4.478"
$ws.Range("E8").Value = "'-1.40%"

# Row 9
$ws.Range("E9").Value = "'-2.18%"

# Row 10
$ws.Range("D10").Value = "'0.9244"

# Row 11
$ws.Range("E11").Value = "'1.86%"

# Row 12
$ws.Range("D12").Value = "'0.1977"
$ws.Range("E12").Value = "'1.56%"

# Row 13
$ws.Range("D13").Value = "'0.09463"
$ws.Range("E13").Value = "'1.26%"

# Row 14
$ws.Range("D14").Value = "'0.03845"
$ws.Range("E14").Value = "'2.38%"

# Row 15
$ws.Range("E15").Value = "'0.88%"

# Row 16
$ws.Range("D16").Value = "'0.001308"
$ws.Range("E16").Value = "'0.48%"

# Row 17
$ws.Range("D17").Value = "'0.006107"
$ws.Range("E17").Value = "'-2.54%"

# Row 18
$ws.Range("D18").Value = "'3.426"
$ws.Range("E18").Value = "'1.90%"

# Row 19
$ws.Range("E19").Value = "'0.67%"

# Row 20
$ws.Range("D20").Value = "'8.850"
$ws.Range("E20").Value = "'1.35%"

# Row 21
$ws.Range("D21").Value = "'0.1363"
$ws.Range("E21").Value = "'-4.10%"

# Row 22
$ws.Range("D22").Value = "'0.2509"
$ws.Range("E22").Value = "'-5.60%"

# Row 23
$ws.Range("D23").Value = "'0.04410"
$ws.Range("E23").Value = "'-0.43%"

# Row 24
$ws.Range("D24").Value = "'0.001271"
$ws.Range("E24").Value = "'0.37%"

# Row 25
$ws.Range("E25").Value = "'1.36%"

# Row 26
$ws.Range("E26").Value = "'-1.84%"

# Row 27
$ws.Range("D27").Value = "'0.0003991"
$ws.Range("E27").Value = "'-0.04%"

# Row 39
$ws.Range("D39").Value = "'0.02870"
$ws.Range("E39").Value = "'0.56%"

# Row 40
$ws.Range("D40").Value = "'0.05523"
$ws.Range("E40").Value = "'0.84%"

# Row 41
$ws.Range("D41").Value = "'0.007969"
$ws.Range("E41").Value = "'2.42%"

# Row 42
$ws.Range("E42").Value = "'1.26%"

# Row 43
$ws.Range("D43").Value = "'0.009002"
$ws.Range("E43").Value = "'-9.84%"

# Row 44
$ws.Range("D44").Value = "'0.002071"
$ws.Range("E44").Value = "'-7.75%"

# Row 45
$ws.Range("D45").Value = "'0.01168"
$ws.Range("E45").Value = "'-0.52%"

# Row 46
$ws.Range("D46").Value = "'0.00006931"
$ws.Range("E46").Value = "'2.29%"

# Row 47
$ws.Range("E47").Value = "'-0.22%"

# Row 48
$ws.Range("D48").Value = "'0.003465"
$ws.Range("E48").Value = "'14.50%"

# Row 49
$ws.Range("D49").Value = "'0.002279"
$ws.Range("E49").Value = "'-0.34%"

# Row 50
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.22%"

# Row 51
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.22%"
